$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.429878711700439
$ws.Range("B1").Value = 2.008875370025635
$ws.Range("C1").Value = 3.066088914871216
$ws.Range("D1").Value = 3.717811346054077
$ws.Range("E1").Value = 0.9659655094146729
